$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Fix the transition text: "retry<3" -> "retry<=3" and drop one dash
#    so the arrow keeps the same overall width, then split the run so a
#    (new) "_GoBack" bookmark sits between "...retry++-------------"
#    and "--> init ".
# ---------------------------------------------------------------------
$full = $d.Content
$foundFull = $full.Find.Execute( `
    "----------retry_init [retry<3] / retry++----------------> init ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundFull) {
    $full.Text = "----------retry_init [retry<=3] / retry++---------------> init "
}

# Split the merged run into "...retry<" | "=3] ... retry++-------------" by
# forcing a run boundary with a temporary bookmark that is immediately removed.
$preEq = $d.Content
$foundPreEq = $preEq.Find.Execute( `
    "----------retry_init [retry<", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundPreEq) {
    $splitPoint = $d.Range($preEq.End, $preEq.End)
    $d.Bookmarks.Add("__TempSplit", $splitPoint) | Out-Null
    $d.Bookmarks("__TempSplit").Delete()
}

# Place the "_GoBack" bookmark right before "--> init " (this also removes
# the bookmark from its previous location automatically, since bookmark
# names are unique in the document).
$preArrow = $d.Content
$foundPreArrow = $preArrow.Find.Execute( `
    "retry++-------------", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundPreArrow) {
    $bmPoint = $d.Range($preArrow.End, $preArrow.End)
    $d.Bookmarks.Add("_GoBack", $bmPoint) | Out-Null
}
